$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text looks numeric must be forced to Text format first,
# otherwise Excel auto-converts the string into a Number and formatting
# like trailing zeros ("97.90" -> 97.9) or ("1.00" -> 1) would be lost.
$ws.Range("D2").Value = "42.912.06"
$ws.Range("E2").Value = "  +0.54%  "
$ws.Range("D3").Value = "2.548.10"
$ws.Range("E3").Value = "  +0.31%  "
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "304.49"
$ws.Range("E5").Value = "  +1.81%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "97.90"
$ws.Range("E6").Value = "  +6.25%  "
$ws.Range("E7").Value = "  +0.67%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("E9").Value = "  -0.36%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.86"
$ws.Range("E10").Value = "  +3.02%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0831"
$ws.Range("E11").Value = "  +3.95%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.75"
$ws.Range("E12").Value = "  +1.81%  "
$ws.Range("E13").Value = "  +1.90%  "
$ws.Range("D14").Value = "2.939.85"
$ws.Range("E14").Value = "  +0.34%  "
$ws.Range("D15").Value = "2.509.87"
$ws.Range("E15").Value = "  -0.33%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.07"
$ws.Range("E16").Value = "  +6.18%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.874"
$ws.Range("E17").Value = "  +0.25%  "
$ws.Range("D18").Value = "42.933.39"
$ws.Range("E18").Value = "  +0.44%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.37"
$ws.Range("E19").Value = "  +3.84%  "
$ws.Range("D20").Value = "0.0₃0994"
$ws.Range("E20").Value = "  +1.41%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.60"
$ws.Range("E21").Value = "  +0.95%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "72.06"
$ws.Range("E22").Value = "  +1.13%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "256.80"
$ws.Range("E23").Value = "  +0.23%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.96"
$ws.Range("E24").Value = "  +1.44%  "
$ws.Range("E25").Value = "  -1.53%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "28.14"
$ws.Range("E26").Value = "  -3.64%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("E27").Value = "  -0.10%  "
$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.29"
$ws.Range("E28").Value = "  +8.44%  "
$ws.Range("B29").Value = "Cosmos"
$ws.Range("C29").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.21"
$ws.Range("E29").Value = "  +1.89%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "37.99"
$ws.Range("E30").Value = "  +2.79%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.14"
$ws.Range("E31").Value = "  +3.58%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "158.45"
$ws.Range("E32").Value = "  +4.01%  "
$ws.Range("E33").Value = "  +14.86%  "
$ws.Range("E34").Value = "  -1.92%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0801"
$ws.Range("E35").Value = "  +1.18%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.32"
$ws.Range("E36").Value = "  -1.84%  "
$ws.Range("E37").Value = "  -4.25%  "
$ws.Range("E38").Value = "  +2.02%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "25.91"
$ws.Range("E39").Value = "  +8.06%  "
$ws.Range("E40").Value = "  +0.34%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.43"
$ws.Range("E41").Value = "  +1.34%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.91"
$ws.Range("E42").Value = "  +1.05%  "
$ws.Range("E43").Value = "  +27.97%  "
$ws.Range("D44").Value = "2.095.44"
$ws.Range("E44").Value = "  +0.84%  "
$ws.Range("E45").Value = "  -1.27%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.00"
$ws.Range("E46").Value = "  +0.09%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "88.15"
$ws.Range("E47").Value = "  +4.56%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.93"
$ws.Range("E48").Value = "  -2.44%  "
$ws.Range("D49").Value = "2.797.31"
$ws.Range("E49").Value = "  +0.34%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "74.88"
$ws.Range("E50").Value = "  +8.90%  "
$ws.Range("E51").Value = "  +1.90%  "
